$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing data rows (2-80), keep header row 1
$ws.Range("A2:G80").ClearContents()

# Write full refreshed dataset for rows 2-85
$ws.Cells.Item(2, 1).Value = 43
$ws.Cells.Item(2, 2).Value = 'How long has Elon Musk been X Corp.''s CEO?'
$ws.Cells.Item(2, 3).Value = 'Elon Musk is no longer X Corp.''s CEO.'
$ws.Cells.Item(2, 4).Value = 'Elon Musk is no longer X Corp.''s CEO.'
$ws.Cells.Item(2, 5).Value = '43.txt'
$ws.Cells.Item(2, 6).Value = '43.txt'
$ws.Cells.Item(2, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(3, 1).Value = 44
$ws.Cells.Item(3, 2).Value = 'Where will the FIFA World Cup be hosted this year?'
$ws.Cells.Item(3, 3).Value = 'There won''t be a FIFA World Cup this year.'
$ws.Cells.Item(3, 4).Value = 'There won''t be a FIFA World Cup this year.'
$ws.Cells.Item(3, 5).Value = '44.txt'
$ws.Cells.Item(3, 6).Value = '44.txt'
$ws.Cells.Item(3, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(4, 1).Value = 92
$ws.Cells.Item(4, 2).Value = 'Alphabet''s market capitalization reached its highest-ever recorded value during what month in 2021?'
$ws.Cells.Item(4, 3).Value = 'The all-time highest value of Alphabet was in April 2024, not in 2021.'
$ws.Cells.Item(4, 4).Value = 'The all-time highest value of Alphabet was in April 2024, not in 2021.'
$ws.Cells.Item(4, 5).Value = '92.txt'
$ws.Cells.Item(4, 6).Value = '92.txt'
$ws.Cells.Item(4, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(5, 1).Value = 95
$ws.Cells.Item(5, 2).Value = 'Which Republican was elected Speaker of the House in January 2023 on the ninth ballot?'
$ws.Cells.Item(5, 3).Value = 'No one received a majority of the votes on the ninth ballot.'
$ws.Cells.Item(5, 4).Value = 'No one received a majority of the votes on the ninth ballot.'
$ws.Cells.Item(5, 5).Value = '95.txt'
$ws.Cells.Item(5, 6).Value = '95.txt'
$ws.Cells.Item(5, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(6, 1).Value = 96
$ws.Cells.Item(6, 2).Value = 'In January 2023, the NHC revised the fatality data of Hurricane Katrina, increasing the reported death toll from 1,800 to what number?'
$ws.Cells.Item(6, 3).Value = 'The reported death toll decreased to 1,392'
$ws.Cells.Item(6, 4).Value = 'The reported death toll decreased to 1,392'
$ws.Cells.Item(6, 5).Value = '96.txt'
$ws.Cells.Item(6, 6).Value = '96.txt'
$ws.Cells.Item(6, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(7, 1).Value = 121
$ws.Cells.Item(7, 2).Value = 'What is the most recent country that President Donald Trump visited during his second presidency?'
$ws.Cells.Item(7, 3).Value = 'President Donald Trump hasn''t visited any countries during his second presidency.'
$ws.Cells.Item(7, 4).Value = 'President Donald Trump hasn''t visited any countries during his second presidency.'
$ws.Cells.Item(7, 5).Value = '121.txt'
$ws.Cells.Item(7, 6).Value = '121.txt'
$ws.Cells.Item(7, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(8, 1).Value = 122
$ws.Cells.Item(8, 2).Value = 'Who was the winner of The Voice US this year?'
$ws.Cells.Item(8, 3).Value = 'This season of The Voice is still ongoing, and the final results have not been announced yet.'
$ws.Cells.Item(8, 4).Value = 'This season of The Voice is still ongoing, and the final results have not been announced yet.'
$ws.Cells.Item(8, 5).Value = '122.txt'
$ws.Cells.Item(8, 6).Value = '122.txt'
$ws.Cells.Item(8, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(9, 1).Value = 123
$ws.Cells.Item(9, 2).Value = 'Who did Michael van Gerwen beat to win this year''s PDC World Darts Championship?'
$ws.Cells.Item(9, 3).Value = 'Michael van Gerwen lost to Luke Littler in the final, held on Friday January 3.'
$ws.Cells.Item(9, 4).Value = 'Michael van Gerwen lost to Luke Littler in the final, held on Friday January 3.'
$ws.Cells.Item(9, 5).Value = '123.txt'
$ws.Cells.Item(9, 6).Value = '123.txt'
$ws.Cells.Item(9, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(10, 1).Value = 124
$ws.Cells.Item(10, 2).Value = 'Who was the winner of American Idol this year?'
$ws.Cells.Item(10, 3).Value = 'This year''s American Idol is still ongoing, and the final results have not been announced yet.'
$ws.Cells.Item(10, 4).Value = 'This year''s American Idol is still ongoing, and the final results have not been announced yet.'
$ws.Cells.Item(10, 5).Value = '124.txt'
$ws.Cells.Item(10, 6).Value = '124.txt'
$ws.Cells.Item(10, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(11, 1).Value = 155
$ws.Cells.Item(11, 2).Value = 'Who is the most recent player to win both the PDC World Youth Championship and the PDC World Darts Championship?'
$ws.Cells.Item(11, 3).Value = 'Luke Littler'
$ws.Cells.Item(11, 4).Value = 'Luke Littler'
$ws.Cells.Item(11, 5).Value = '155.txt'
$ws.Cells.Item(11, 6).Value = '155.txt'
$ws.Cells.Item(11, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(12, 1).Value = 156
$ws.Cells.Item(12, 2).Value = 'What was the Weeknd''s last studio album titled?'
$ws.Cells.Item(12, 3).Value = 'Hurry Up Tomorrow'
$ws.Cells.Item(12, 4).Value = 'Hurry Up Tomorrow'
$ws.Cells.Item(12, 5).Value = '156.txt'
$ws.Cells.Item(12, 6).Value = '156.txt'
$ws.Cells.Item(12, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(13, 1).Value = 157
$ws.Cells.Item(13, 2).Value = 'How many seats are there in the German Bundestag?'
$ws.Cells.Item(13, 3).NumberFormat = "@"
$ws.Cells.Item(13, 3).Value = '630'
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '630'
$ws.Cells.Item(13, 5).Value = '157.txt'
$ws.Cells.Item(13, 6).Value = '157.txt'
$ws.Cells.Item(13, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(14, 1).Value = 160
$ws.Cells.Item(14, 2).Value = 'Who is the President of the United States?'
$ws.Cells.Item(14, 3).Value = 'Donald Trump'
$ws.Cells.Item(14, 4).Value = 'Donald Trump'
$ws.Cells.Item(14, 5).Value = '160.txt'
$ws.Cells.Item(14, 6).Value = '160.txt'
$ws.Cells.Item(14, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(15, 1).Value = 162
$ws.Cells.Item(15, 2).Value = 'Who is the current Prime Minister of Japan?'
$ws.Cells.Item(15, 3).Value = 'Shigeru Ishiba'
$ws.Cells.Item(15, 4).Value = 'Shigeru Ishiba'
$ws.Cells.Item(15, 5).Value = '162.txt'
$ws.Cells.Item(15, 6).Value = '162.txt'
$ws.Cells.Item(15, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(16, 1).Value = 163
$ws.Cells.Item(16, 2).Value = 'How many asteroids have been discovered before impacting Earth?'
$ws.Cells.Item(16, 3).NumberFormat = "@"
$ws.Cells.Item(16, 3).Value = '11'
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '11'
$ws.Cells.Item(16, 5).Value = '163.txt'
$ws.Cells.Item(16, 6).Value = '163.txt'
$ws.Cells.Item(16, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(17, 1).Value = 164
$ws.Cells.Item(17, 2).Value = 'Which country is the most recent member state of the Schengen Area?'
$ws.Cells.Item(17, 3).Value = 'Bulgaria and Romania'
$ws.Cells.Item(17, 4).Value = 'Bulgaria and Romania'
$ws.Cells.Item(17, 5).Value = '164.txt'
$ws.Cells.Item(17, 6).Value = '164.txt'
$ws.Cells.Item(17, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(18, 1).Value = 165
$ws.Cells.Item(18, 2).Value = 'What is Greg Egan''s latest novel?'
$ws.Cells.Item(18, 3).Value = 'Morphotrophic'
$ws.Cells.Item(18, 4).Value = 'Morphotrophic'
$ws.Cells.Item(18, 5).Value = '165.txt'
$ws.Cells.Item(18, 6).Value = '165.txt'
$ws.Cells.Item(18, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(19, 1).Value = 166
$ws.Cells.Item(19, 2).Value = 'What is the title of Taylor Swift''s most recent studio album?'
$ws.Cells.Item(19, 3).Value = 'The Tortured Poets Department'
$ws.Cells.Item(19, 4).Value = 'The Tortured Poets Department'
$ws.Cells.Item(19, 5).Value = '166.txt'
$ws.Cells.Item(19, 6).Value = '166.txt'
$ws.Cells.Item(19, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(20, 1).Value = 167
$ws.Cells.Item(20, 2).Value = 'Who is the current Prime Minister of the United Kingdom?'
$ws.Cells.Item(20, 3).Value = 'Keir Starmer'
$ws.Cells.Item(20, 4).Value = 'Keir Starmer'
$ws.Cells.Item(20, 5).Value = '167.txt'
$ws.Cells.Item(20, 6).Value = '167.txt'
$ws.Cells.Item(20, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(21, 1).Value = 168
$ws.Cells.Item(21, 2).Value = 'Who is the most recent former Prime Minister of the United Kingdom?'
$ws.Cells.Item(21, 3).Value = 'Rishi Sunak'
$ws.Cells.Item(21, 4).Value = 'Rishi Sunak'
$ws.Cells.Item(21, 5).Value = '168.txt'
$ws.Cells.Item(21, 6).Value = '168.txt'
$ws.Cells.Item(21, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(22, 1).Value = 170
$ws.Cells.Item(22, 2).Value = 'What is the most recent outbreak declared as a public health emergency of international concern by WHO?'
$ws.Cells.Item(22, 3).Value = 'Clade I mpox'
$ws.Cells.Item(22, 4).Value = 'Clade I mpox'
$ws.Cells.Item(22, 5).Value = '170.txt'
$ws.Cells.Item(22, 6).Value = '170.txt'
$ws.Cells.Item(22, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(23, 1).Value = 172
$ws.Cells.Item(23, 2).Value = 'How many LA metro rail stations are there?'
$ws.Cells.Item(23, 3).NumberFormat = "@"
$ws.Cells.Item(23, 3).Value = '102'
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '102'
$ws.Cells.Item(23, 5).Value = '172.txt'
$ws.Cells.Item(23, 6).Value = '172.txt'
$ws.Cells.Item(23, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(24, 1).Value = 173
$ws.Cells.Item(24, 2).Value = 'Who won the biggest single-tournament payday in tennis history?'
$ws.Cells.Item(24, 3).Value = 'Jannik Sinner'
$ws.Cells.Item(24, 4).Value = 'Jannik Sinner'
$ws.Cells.Item(24, 5).Value = '173.txt'
$ws.Cells.Item(24, 6).Value = '173.txt'
$ws.Cells.Item(24, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(25, 1).Value = 175
$ws.Cells.Item(25, 2).Value = 'On what date did Kendrick Lamar release his most recent studio album?'
$ws.Cells.Item(25, 3).Value = 'Nov 22, 2024'
$ws.Cells.Item(25, 4).Value = 'Nov 22, 2024'
$ws.Cells.Item(25, 5).Value = '175.txt'
$ws.Cells.Item(25, 6).Value = '175.txt'
$ws.Cells.Item(25, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(26, 1).Value = 180
$ws.Cells.Item(26, 2).Value = 'What''s the latest major version of the .NET?'
$ws.Cells.Item(26, 3).Value = '.NET 9'
$ws.Cells.Item(26, 4).Value = '.NET 9'
$ws.Cells.Item(26, 5).Value = '180.txt'
$ws.Cells.Item(26, 6).Value = '180.txt'
$ws.Cells.Item(26, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(27, 1).Value = 182
$ws.Cells.Item(27, 2).Value = 'How many food allergens with mandatory labeling are there in the United States?'
$ws.Cells.Item(27, 3).NumberFormat = "@"
$ws.Cells.Item(27, 3).Value = '9'
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '9'
$ws.Cells.Item(27, 5).Value = '182.txt'
$ws.Cells.Item(27, 6).Value = '182.txt'
$ws.Cells.Item(27, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(28, 1).Value = 183
$ws.Cells.Item(28, 2).Value = 'Who is the X Corp.''s CEO?'
$ws.Cells.Item(28, 3).Value = 'Linda Yaccarino'
$ws.Cells.Item(28, 4).Value = 'Linda Yaccarino'
$ws.Cells.Item(28, 5).Value = '183.txt'
$ws.Cells.Item(28, 6).Value = '183.txt'
$ws.Cells.Item(28, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(29, 1).Value = 184
$ws.Cells.Item(29, 2).Value = 'Who has the most followers on Twitter?'
$ws.Cells.Item(29, 3).Value = 'Elon Musk'
$ws.Cells.Item(29, 4).Value = 'Elon Musk'
$ws.Cells.Item(29, 5).Value = '184.txt'
$ws.Cells.Item(29, 6).Value = '184.txt'
$ws.Cells.Item(29, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(30, 1).Value = 185
$ws.Cells.Item(30, 2).Value = 'What''s the largest stadium by capacity in the world?'
$ws.Cells.Item(30, 3).Value = 'Narendra Modi Stadium'
$ws.Cells.Item(30, 4).Value = 'Narendra Modi Stadium'
$ws.Cells.Item(30, 5).Value = '185.txt'
$ws.Cells.Item(30, 6).Value = '185.txt'
$ws.Cells.Item(30, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(31, 1).Value = 187
$ws.Cells.Item(31, 2).Value = 'Who is the chancellor of UMass Amherst?'
$ws.Cells.Item(31, 3).Value = 'Javier Reyes'
$ws.Cells.Item(31, 4).Value = 'Javier Reyes'
$ws.Cells.Item(31, 5).Value = '187.txt'
$ws.Cells.Item(31, 6).Value = '187.txt'
$ws.Cells.Item(31, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(32, 1).Value = 189
$ws.Cells.Item(32, 2).Value = 'What is the most recently released Studio Ghibli film?'
$ws.Cells.Item(32, 3).Value = 'The Boy and the Heron'
$ws.Cells.Item(32, 4).Value = 'The Boy and the Heron'
$ws.Cells.Item(32, 5).Value = '189.txt'
$ws.Cells.Item(32, 6).Value = '189.txt'
$ws.Cells.Item(32, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(33, 1).Value = 192
$ws.Cells.Item(33, 2).Value = 'What is Croatia''s current national currency?'
$ws.Cells.Item(33, 3).Value = 'Euro'
$ws.Cells.Item(33, 4).Value = 'Euro'
$ws.Cells.Item(33, 5).Value = '192.txt'
$ws.Cells.Item(33, 6).Value = '192.txt'
$ws.Cells.Item(33, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(34, 1).Value = 193
$ws.Cells.Item(34, 2).Value = 'How many member states are there in the Eurozone?'
$ws.Cells.Item(34, 3).NumberFormat = "@"
$ws.Cells.Item(34, 3).Value = '20'
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '20'
$ws.Cells.Item(34, 5).Value = '193.txt'
$ws.Cells.Item(34, 6).Value = '193.txt'
$ws.Cells.Item(34, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(35, 1).Value = 196
$ws.Cells.Item(35, 2).Value = 'What country became the latest country to adopt the Euro?'
$ws.Cells.Item(35, 3).Value = 'Croatia'
$ws.Cells.Item(35, 4).Value = 'Croatia'
$ws.Cells.Item(35, 5).Value = '196.txt'
$ws.Cells.Item(35, 6).Value = '196.txt'
$ws.Cells.Item(35, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(36, 1).Value = 198
$ws.Cells.Item(36, 2).Value = 'How many vehicle models does Tesla offer?'
$ws.Cells.Item(36, 3).NumberFormat = "@"
$ws.Cells.Item(36, 3).Value = '6'
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '6'
$ws.Cells.Item(36, 5).Value = '198.txt'
$ws.Cells.Item(36, 6).Value = '198.txt'
$ws.Cells.Item(36, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(37, 1).Value = 199
$ws.Cells.Item(37, 2).Value = 'Which country won the latest World Pastry Cup?'
$ws.Cells.Item(37, 3).Value = 'Japan'
$ws.Cells.Item(37, 4).Value = 'Japan'
$ws.Cells.Item(37, 5).Value = '199.txt'
$ws.Cells.Item(37, 6).Value = '199.txt'
$ws.Cells.Item(37, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(38, 1).Value = 377
$ws.Cells.Item(38, 2).Value = 'How many exoplanets have been discovered?'
$ws.Cells.Item(38, 3).Value = '5,867 confirmed exoplanets'
$ws.Cells.Item(38, 4).Value = '5,867 confirmed exoplanets'
$ws.Cells.Item(38, 5).Value = '377.txt'
$ws.Cells.Item(38, 6).Value = '377.txt'
$ws.Cells.Item(38, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(39, 1).Value = 380
$ws.Cells.Item(39, 2).Value = 'How many tornadoes have been confirmed by Enhanced Fujita rating in the United States so far this year?'
$ws.Cells.Item(39, 3).NumberFormat = "@"
$ws.Cells.Item(39, 3).Value = '193'
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '193'
$ws.Cells.Item(39, 5).Value = '380.txt'
$ws.Cells.Item(39, 6).Value = '380.txt'
$ws.Cells.Item(39, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(40, 1).Value = 384
$ws.Cells.Item(40, 2).Value = 'What is the current most popular Wikipedia article of the Week?'
$ws.Cells.Item(40, 3).Value = 'Adolescence (TV series)'
$ws.Cells.Item(40, 4).Value = 'Adolescence (TV series)'
$ws.Cells.Item(40, 5).Value = '384.txt'
$ws.Cells.Item(40, 6).Value = '384.txt'
$ws.Cells.Item(40, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(41, 1).Value = 391
$ws.Cells.Item(41, 2).Value = 'Who won the most recent Super Bowl?'
$ws.Cells.Item(41, 3).Value = 'Philadelphia Eagles'
$ws.Cells.Item(41, 4).Value = 'Philadelphia Eagles'
$ws.Cells.Item(41, 5).Value = '391.txt'
$ws.Cells.Item(41, 6).Value = '391.txt'
$ws.Cells.Item(41, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(42, 1).Value = 393
$ws.Cells.Item(42, 2).Value = 'What are this year''s Coachella festival dates?'
$ws.Cells.Item(42, 3).Value = 'April 11–13 and April 18–20'
$ws.Cells.Item(42, 4).Value = 'April 11–13 and April 18–20'
$ws.Cells.Item(42, 5).Value = '393.txt'
$ws.Cells.Item(42, 6).Value = '393.txt'
$ws.Cells.Item(42, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(43, 1).Value = 401
$ws.Cells.Item(43, 2).Value = 'Which broadcast network is this year''s Super Bowl broadcaster?'
$ws.Cells.Item(43, 3).Value = 'Fox and Tubi'
$ws.Cells.Item(43, 4).Value = 'Fox and Tubi'
$ws.Cells.Item(43, 5).Value = '401.txt'
$ws.Cells.Item(43, 6).Value = '401.txt'
$ws.Cells.Item(43, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(44, 1).Value = 402
$ws.Cells.Item(44, 2).Value = 'What was the highest-grossing film in the United States last year?'
$ws.Cells.Item(44, 3).Value = 'Inside Out 2'
$ws.Cells.Item(44, 4).Value = 'Inside Out 2'
$ws.Cells.Item(44, 5).Value = '402.txt'
$ws.Cells.Item(44, 6).Value = '402.txt'
$ws.Cells.Item(44, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(45, 1).Value = 403
$ws.Cells.Item(45, 2).Value = 'What Berber year corresponds to the present year?'
$ws.Cells.Item(45, 3).NumberFormat = "@"
$ws.Cells.Item(45, 3).Value = '2975'
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '2975'
$ws.Cells.Item(45, 5).Value = '403.txt'
$ws.Cells.Item(45, 6).Value = '403.txt'
$ws.Cells.Item(45, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(46, 1).Value = 408
$ws.Cells.Item(46, 2).Value = 'On what date was last year''s State of the Union Address delivered?'
$ws.Cells.Item(46, 3).Value = 'Mar 7, 2024'
$ws.Cells.Item(46, 4).Value = 'Mar 7, 2024'
$ws.Cells.Item(46, 5).Value = '408.txt'
$ws.Cells.Item(46, 6).Value = '408.txt'
$ws.Cells.Item(46, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(47, 1).Value = 409
$ws.Cells.Item(47, 2).Value = 'Who won the latest Grammy Award for Best Jazz Performance?'
$ws.Cells.Item(47, 3).Value = 'Samara Joy featuring Sullivan Fortner'
$ws.Cells.Item(47, 4).Value = 'Samara Joy featuring Sullivan Fortner'
$ws.Cells.Item(47, 5).Value = '409.txt'
$ws.Cells.Item(47, 6).Value = '409.txt'
$ws.Cells.Item(47, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(48, 1).Value = 410
$ws.Cells.Item(48, 2).Value = 'What won the latest Academy Award for Best Animated Feature?'
$ws.Cells.Item(48, 3).Value = 'Flow'
$ws.Cells.Item(48, 4).Value = 'Flow'
$ws.Cells.Item(48, 5).Value = '410.txt'
$ws.Cells.Item(48, 6).Value = '410.txt'
$ws.Cells.Item(48, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(49, 1).Value = 411
$ws.Cells.Item(49, 2).Value = 'Which movie won the latest Academy Award for Best Picture?'
$ws.Cells.Item(49, 3).Value = 'Anora'
$ws.Cells.Item(49, 4).Value = 'Anora'
$ws.Cells.Item(49, 5).Value = '411.txt'
$ws.Cells.Item(49, 6).Value = '411.txt'
$ws.Cells.Item(49, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(50, 1).Value = 420
$ws.Cells.Item(50, 2).Value = 'Which team is at the top of the latest Premier League season?'
$ws.Cells.Item(50, 3).Value = 'Liverpool'
$ws.Cells.Item(50, 4).Value = 'Liverpool'
$ws.Cells.Item(50, 5).Value = '420.txt'
$ws.Cells.Item(50, 6).Value = '420.txt'
$ws.Cells.Item(50, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(51, 1).Value = 421
$ws.Cells.Item(51, 2).Value = 'Is Arsenal on the top of the latest Premier League standings?'
$ws.Cells.Item(51, 3).Value = 'No'
$ws.Cells.Item(51, 4).Value = 'No'
$ws.Cells.Item(51, 5).Value = '421.txt'
$ws.Cells.Item(51, 6).Value = '421.txt'
$ws.Cells.Item(51, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(52, 1).Value = 422
$ws.Cells.Item(52, 2).Value = 'Who was the guest of honor at the most recent state dinner hosted by the President of the United States?'
$ws.Cells.Item(52, 3).Value = 'William Ruto'
$ws.Cells.Item(52, 4).Value = 'William Ruto'
$ws.Cells.Item(52, 5).Value = '422.txt'
$ws.Cells.Item(52, 6).Value = '422.txt'
$ws.Cells.Item(52, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(53, 1).Value = 423
$ws.Cells.Item(53, 2).Value = 'What album holds the record for the most Spotify streams reached in a single day?'
$ws.Cells.Item(53, 3).Value = 'The Tortured Poets Department'
$ws.Cells.Item(53, 4).Value = 'The Tortured Poets Department'
$ws.Cells.Item(53, 5).Value = '423.txt'
$ws.Cells.Item(53, 6).Value = '423.txt'
$ws.Cells.Item(53, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(54, 1).Value = 424
$ws.Cells.Item(54, 2).Value = 'Who was the most recent incumbent president worldwide who ran for re-election but was not reelected?'
$ws.Cells.Item(54, 3).Value = 'Muse Bihi Abdi'
$ws.Cells.Item(54, 4).Value = 'Muse Bihi Abdi'
$ws.Cells.Item(54, 5).Value = '424.txt'
$ws.Cells.Item(54, 6).Value = '424.txt'
$ws.Cells.Item(54, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(55, 1).Value = 427
$ws.Cells.Item(55, 2).Value = 'Who is the current ATP top-ranked men''s singles tennis player?'
$ws.Cells.Item(55, 3).Value = 'Jannik Sinner'
$ws.Cells.Item(55, 4).Value = 'Jannik Sinner'
$ws.Cells.Item(55, 5).Value = '427.txt'
$ws.Cells.Item(55, 6).Value = '427.txt'
$ws.Cells.Item(55, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(56, 1).Value = 430
$ws.Cells.Item(56, 2).Value = 'Who is the number 1 ranked female tennis player in the world?'
$ws.Cells.Item(56, 3).Value = 'Aryna Sabalenka'
$ws.Cells.Item(56, 4).Value = 'Aryna Sabalenka'
$ws.Cells.Item(56, 5).Value = '430.txt'
$ws.Cells.Item(56, 6).Value = '430.txt'
$ws.Cells.Item(56, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(57, 1).Value = 432
$ws.Cells.Item(57, 2).Value = 'Who is the latest MotoGP World Riders'' Champion?'
$ws.Cells.Item(57, 3).Value = 'Jorge Martín'
$ws.Cells.Item(57, 4).Value = 'Jorge Martín'
$ws.Cells.Item(57, 5).Value = '432.txt'
$ws.Cells.Item(57, 6).Value = '432.txt'
$ws.Cells.Item(57, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(58, 1).Value = 433
$ws.Cells.Item(58, 2).Value = 'Who won the most recent Time Magazine''s Athlete of the Year?'
$ws.Cells.Item(58, 3).Value = 'Caitlin Clark'
$ws.Cells.Item(58, 4).Value = 'Caitlin Clark'
$ws.Cells.Item(58, 5).Value = '433.txt'
$ws.Cells.Item(58, 6).Value = '433.txt'
$ws.Cells.Item(58, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(59, 1).Value = 435
$ws.Cells.Item(59, 2).Value = 'What book won the latest Nebula award for Best Novel?'
$ws.Cells.Item(59, 3).Value = 'The Saint of Bright Doors'
$ws.Cells.Item(59, 4).Value = 'The Saint of Bright Doors'
$ws.Cells.Item(59, 5).Value = '435.txt'
$ws.Cells.Item(59, 6).Value = '435.txt'
$ws.Cells.Item(59, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(60, 1).Value = 436
$ws.Cells.Item(60, 2).Value = 'Which game won the Spiel des Jahres award most recently?'
$ws.Cells.Item(60, 3).Value = 'Sky Team'
$ws.Cells.Item(60, 4).Value = 'Sky Team'
$ws.Cells.Item(60, 5).Value = '436.txt'
$ws.Cells.Item(60, 6).Value = '436.txt'
$ws.Cells.Item(60, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(61, 1).Value = 437
$ws.Cells.Item(61, 2).Value = 'Which drama series won the most recent Primetime Emmy Award for Outstanding Drama Series?'
$ws.Cells.Item(61, 3).Value = 'Shōgun'
$ws.Cells.Item(61, 4).Value = 'Shōgun'
$ws.Cells.Item(61, 5).Value = '437.txt'
$ws.Cells.Item(61, 6).Value = '437.txt'
$ws.Cells.Item(61, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(62, 1).Value = 438
$ws.Cells.Item(62, 2).Value = 'What is the name of the most recent episode of The Lord of the Rings: The Rings of Power?'
$ws.Cells.Item(62, 3).Value = 'Shadow and Flame'
$ws.Cells.Item(62, 4).Value = 'Shadow and Flame'
$ws.Cells.Item(62, 5).Value = '438.txt'
$ws.Cells.Item(62, 6).Value = '438.txt'
$ws.Cells.Item(62, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(63, 1).Value = 439
$ws.Cells.Item(63, 2).Value = 'On what date did the Patriots last play the Miami Dolphins?'
$ws.Cells.Item(63, 3).Value = 'Nov 24, 2024'
$ws.Cells.Item(63, 4).Value = 'Nov 24, 2024'
$ws.Cells.Item(63, 5).Value = '439.txt'
$ws.Cells.Item(63, 6).Value = '439.txt'
$ws.Cells.Item(63, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(64, 1).Value = 441
$ws.Cells.Item(64, 2).Value = 'What is the next Walt Disney Animation Studios film?'
$ws.Cells.Item(64, 3).Value = 'Zootopia 2'
$ws.Cells.Item(64, 4).Value = 'Zootopia 2'
$ws.Cells.Item(64, 5).Value = '441.txt'
$ws.Cells.Item(64, 6).Value = '441.txt'
$ws.Cells.Item(64, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(65, 1).Value = 442
$ws.Cells.Item(65, 2).Value = 'What is the most recent Walt Disney Animation Studios''s animated film?'
$ws.Cells.Item(65, 3).Value = 'Moana 2'
$ws.Cells.Item(65, 4).Value = 'Moana 2'
$ws.Cells.Item(65, 5).Value = '442.txt'
$ws.Cells.Item(65, 6).Value = '442.txt'
$ws.Cells.Item(65, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(66, 1).Value = 444
$ws.Cells.Item(66, 2).Value = 'When did the latest NFL season begin?'
$ws.Cells.Item(66, 3).Value = 'Sep 5, 2024'
$ws.Cells.Item(66, 4).Value = 'Sep 5, 2024'
$ws.Cells.Item(66, 5).Value = '444.txt'
$ws.Cells.Item(66, 6).Value = '444.txt'
$ws.Cells.Item(66, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(67, 1).Value = 446
$ws.Cells.Item(67, 2).Value = 'Where will the upcoming International Mathematical Olympiad (IMO) be hosted?'
$ws.Cells.Item(67, 3).Value = 'Sunshine Coast'
$ws.Cells.Item(67, 4).Value = 'Sunshine Coast'
$ws.Cells.Item(67, 5).Value = '446.txt'
$ws.Cells.Item(67, 6).Value = '446.txt'
$ws.Cells.Item(67, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(68, 1).Value = 448
$ws.Cells.Item(68, 2).Value = 'Who won the latest NBA championship?'
$ws.Cells.Item(68, 3).Value = 'Boston Celtics'
$ws.Cells.Item(68, 4).Value = 'Boston Celtics'
$ws.Cells.Item(68, 5).Value = '448.txt'
$ws.Cells.Item(68, 6).Value = '448.txt'
$ws.Cells.Item(68, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(69, 1).Value = 449
$ws.Cells.Item(69, 2).Value = 'The longest winning streak in the last season of the Premier League consisted of how many matches?'
$ws.Cells.Item(69, 3).NumberFormat = "@"
$ws.Cells.Item(69, 3).Value = '9'
$ws.Cells.Item(69, 4).NumberFormat = "@"
$ws.Cells.Item(69, 4).Value = '9'
$ws.Cells.Item(69, 5).Value = '449.txt'
$ws.Cells.Item(69, 6).Value = '449.txt'
$ws.Cells.Item(69, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(70, 1).Value = 452
$ws.Cells.Item(70, 2).Value = 'How many books has Colleen Hoover published?'
$ws.Cells.Item(70, 3).Value = '26 books'
$ws.Cells.Item(70, 4).Value = '26 books'
$ws.Cells.Item(70, 5).Value = '452.txt'
$ws.Cells.Item(70, 6).Value = '452.txt'
$ws.Cells.Item(70, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(71, 1).Value = 453
$ws.Cells.Item(71, 2).Value = 'What is the latest United States jurisdiction to legalize the recreational use of cannabis?'
$ws.Cells.Item(71, 3).Value = 'Ohio'
$ws.Cells.Item(71, 4).Value = 'Ohio'
$ws.Cells.Item(71, 5).Value = '453.txt'
$ws.Cells.Item(71, 6).Value = '453.txt'
$ws.Cells.Item(71, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(72, 1).Value = 532
$ws.Cells.Item(72, 2).Value = 'What is the hottest year on record?'
$ws.Cells.Item(72, 3).NumberFormat = "@"
$ws.Cells.Item(72, 3).Value = '2024'
$ws.Cells.Item(72, 4).NumberFormat = "@"
$ws.Cells.Item(72, 4).Value = '2024'
$ws.Cells.Item(72, 5).Value = '532.txt'
$ws.Cells.Item(72, 6).Value = '532.txt'
$ws.Cells.Item(72, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(73, 1).Value = 535
$ws.Cells.Item(73, 2).Value = 'Who is the current Jeopardy! host?'
$ws.Cells.Item(73, 3).Value = 'Ken Jennings'
$ws.Cells.Item(73, 4).Value = 'Ken Jennings'
$ws.Cells.Item(73, 5).Value = '535.txt'
$ws.Cells.Item(73, 6).Value = '535.txt'
$ws.Cells.Item(73, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(74, 1).Value = 536
$ws.Cells.Item(74, 2).Value = 'Who holds the record for most hundreds in one day international men''s cricket?'
$ws.Cells.Item(74, 3).Value = 'Virat Kohli'
$ws.Cells.Item(74, 4).Value = 'Virat Kohli'
$ws.Cells.Item(74, 5).Value = '536.txt'
$ws.Cells.Item(74, 6).Value = '536.txt'
$ws.Cells.Item(74, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(75, 1).Value = 539
$ws.Cells.Item(75, 2).Value = 'What was the largest acquisition deal of a tech company in history?'
$ws.Cells.Item(75, 3).Value = 'Microsoft purchases Activision Blizzard for US$68.7 billion.'
$ws.Cells.Item(75, 4).Value = 'Microsoft purchases Activision Blizzard for US$68.7 billion.'
$ws.Cells.Item(75, 5).Value = '539.txt'
$ws.Cells.Item(75, 6).Value = '539.txt'
$ws.Cells.Item(75, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(76, 1).Value = 576
$ws.Cells.Item(76, 2).Value = 'When is the next Lunar New Year?'
$ws.Cells.Item(76, 3).Value = 'Feb 17, 2026'
$ws.Cells.Item(76, 4).Value = 'Feb 17, 2026'
$ws.Cells.Item(76, 5).Value = '576.txt'
$ws.Cells.Item(76, 6).Value = '576.txt'
$ws.Cells.Item(76, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(77, 1).Value = 577
$ws.Cells.Item(77, 2).Value = 'When does this year''s Rosh Hashanah start?'
$ws.Cells.Item(77, 3).Value = 'Sep 22, 2026'
$ws.Cells.Item(77, 4).Value = 'Sep 22, 2026'
$ws.Cells.Item(77, 5).Value = '577.txt'
$ws.Cells.Item(77, 6).Value = '577.txt'
$ws.Cells.Item(77, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(78, 1).Value = 578
$ws.Cells.Item(78, 2).Value = 'What is the host city of the upcoming International Olympiad in Informatics (IOI)?'
$ws.Cells.Item(78, 3).Value = 'Sucre'
$ws.Cells.Item(78, 4).Value = 'Sucre'
$ws.Cells.Item(78, 5).Value = '578.txt'
$ws.Cells.Item(78, 6).Value = '578.txt'
$ws.Cells.Item(78, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(79, 1).Value = 579
$ws.Cells.Item(79, 2).Value = 'What won the latest Academy Award for Best Animated Short Film?'
$ws.Cells.Item(79, 3).Value = 'In the Shadow of the Cypress'
$ws.Cells.Item(79, 4).Value = 'In the Shadow of the Cypress'
$ws.Cells.Item(79, 5).Value = '579.txt'
$ws.Cells.Item(79, 6).Value = '579.txt'
$ws.Cells.Item(79, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(80, 1).Value = 580
$ws.Cells.Item(80, 2).Value = 'What is the next Columbia Pictures film?'
$ws.Cells.Item(80, 3).Value = 'Karate Kid: Legends'
$ws.Cells.Item(80, 4).Value = 'Karate Kid: Legends'
$ws.Cells.Item(80, 5).Value = '580.txt'
$ws.Cells.Item(80, 6).Value = '580.txt'
$ws.Cells.Item(80, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(81, 1).Value = 586
$ws.Cells.Item(81, 2).Value = 'What is the name of the most recent hurricane that affected the Southeastern Coast of the United States?'
$ws.Cells.Item(81, 3).Value = 'Milton'
$ws.Cells.Item(81, 4).Value = 'Milton'
$ws.Cells.Item(81, 5).Value = '586.txt'
$ws.Cells.Item(81, 6).Value = '586.txt'
$ws.Cells.Item(81, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(82, 1).Value = 587
$ws.Cells.Item(82, 2).Value = 'What is King Gizzard’s most recent studio album?'
$ws.Cells.Item(82, 3).Value = 'Flight b741'
$ws.Cells.Item(82, 4).Value = 'Flight b741'
$ws.Cells.Item(82, 5).Value = '587.txt'
$ws.Cells.Item(82, 6).Value = '587.txt'
$ws.Cells.Item(82, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(83, 1).Value = 588
$ws.Cells.Item(83, 2).Value = 'Which comedy series won the most recent Primetime Emmy Award for Outstanding Comedy Series?'
$ws.Cells.Item(83, 3).Value = 'Hacks (Season 3)'
$ws.Cells.Item(83, 4).Value = 'Hacks (Season 3)'
$ws.Cells.Item(83, 5).Value = '588.txt'
$ws.Cells.Item(83, 6).Value = '588.txt'
$ws.Cells.Item(83, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(84, 1).Value = 589
$ws.Cells.Item(84, 2).Value = 'What institution won the most recent ACM-ICPC World Finals?'
$ws.Cells.Item(84, 3).Value = 'Peking University'
$ws.Cells.Item(84, 4).Value = 'Peking University'
$ws.Cells.Item(84, 5).Value = '589.txt'
$ws.Cells.Item(84, 6).Value = '589.txt'
$ws.Cells.Item(84, 7).Value = '[''test1'', ''test2'', ''test3'']'
$ws.Cells.Item(85, 1).Value = 590
$ws.Cells.Item(85, 2).Value = 'Who won the most recent season of America''s Got Talent?'
$ws.Cells.Item(85, 3).Value = 'Richard Goodall'
$ws.Cells.Item(85, 4).Value = 'Richard Goodall'
$ws.Cells.Item(85, 5).Value = '590.txt'
$ws.Cells.Item(85, 6).Value = '590.txt'
$ws.Cells.Item(85, 7).Value = '[''test1'', ''test2'', ''test3'']'

Write-Host "Done. Final row count (excluding header): " (84)